$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.524.49"
$ws.Range("E2").Value = "'  -0.71%  "
$ws.Range("D3").Value = "'1.850.60"
$ws.Range("E3").Value = "'  -0.33%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("D5").Value = "'243.03"
$ws.Range("E5").Value = "'  -0.74%  "
$ws.Range("D6").Value = "'0.6356"
$ws.Range("E6").Value = "'  -1.39%  "
$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("D8").Value = "'0.07560"
$ws.Range("E8").Value = "'  +0.84%  "
$ws.Range("D9").Value = "'0.2974"
$ws.Range("E9").Value = "'  -0.27%  "
$ws.Range("D10").Value = "'24.28"
$ws.Range("E10").Value = "'  -0.53%  "
$ws.Range("D11").Value = "'0.07689"
$ws.Range("E11").Value = "'  +0.40%  "
$ws.Range("D12").Value = "'1.857.21"
$ws.Range("E12").Value = "'  +0.02%  "
$ws.Range("D13").Value = "'5.026"
$ws.Range("E13").Value = "'  -0.32%  "
$ws.Range("D14").Value = "'0.6869"
$ws.Range("E14").Value = "'  -0.55%  "
$ws.Range("D15").Value = "'83.93"
$ws.Range("E15").Value = "'  +0.00%  "
$ws.Range("D16").Value = "'0.000009773"
$ws.Range("E16").Value = "'  +1.48%  "
$ws.Range("D17").Value = "'2.117.57"
$ws.Range("E17").Value = "'  +0.25%  "
$ws.Range("D18").Value = "'6.221"
$ws.Range("D19").Value = "'29.543.15"
$ws.Range("E19").Value = "'  -0.61%  "
$ws.Range("D20").Value = "'236.05"
$ws.Range("E20").Value = "'  +0.11%  "
$ws.Range("D21").Value = "'12.53"
$ws.Range("E21").Value = "'  -1.08%  "
$ws.Range("D23").Value = "'7.639"
$ws.Range("E23").Value = "'  +2.41%  "
$ws.Range("D24").Value = "'0.9998"
$ws.Range("E24").Value = "'  +0.06%  "
$ws.Range("D25").Value = "'156.06"
$ws.Range("E25").Value = "'  -1.59%  "
$ws.Range("D26").Value = "'0.1392"
$ws.Range("E26").Value = "'  -1.91%  "
$ws.Range("D27").Value = "'8.463"
$ws.Range("E27").Value = "'  -0.89%  "
$ws.Range("D28").Value = "'17.76"
$ws.Range("E28").Value = "'  -0.87%  "
$ws.Range("D29").Value = "'1.486"
$ws.Range("E29").Value = "'  -0.57%  "
$ws.Range("D30").Value = "'0.05849"
$ws.Range("E30").Value = "'  -7.08%  "
$ws.Range("D31").Value = "'1.277"
$ws.Range("E31").Value = "'  -0.09%  "
$ws.Range("D32").Value = "'4.126"
$ws.Range("E32").Value = "'  -0.61%  "
$ws.Range("D33").Value = "'4.050"
$ws.Range("E33").Value = "'  -0.85%  "
$ws.Range("D34").Value = "'1.900"
$ws.Range("E34").Value = "'  +0.12%  "
$ws.Range("D35").Value = "'1.173"
$ws.Range("E35").Value = "'  +0.05%  "
$ws.Range("D36").Value = "'0.7188"
$ws.Range("E36").Value = "'  -1.48%  "
$ws.Range("D37").Value = "'2.597"
$ws.Range("E37").Value = "'  -0.53%  "
$ws.Range("D38").Value = "'1.243.66"
$ws.Range("E38").Value = "'  +3.38%  "
$ws.Range("D39").Value = "'2.802"
$ws.Range("E39").Value = "'  -1.38%  "
$ws.Range("D40").Value = "'0.01775"
$ws.Range("E40").Value = "'  -0.69%  "
$ws.Range("D41").Value = "'0.9135"
$ws.Range("E41").Value = "'  -0.95%  "
$ws.Range("D42").Value = "'6.119"
$ws.Range("E42").Value = "'  -0.36%  "
$ws.Range("B43").Value = "'PaxDollar"
$ws.Range("C43").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'0.9994"
$ws.Range("E43").Value = "'  -0.04%  "
$ws.Range("B44").Value = "'RocketPoolETH"
$ws.Range("C44").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "'2.026.78"
$ws.Range("E44").Value = "'  +0.30%  "
$ws.Range("D45").Value = "'102.46"
$ws.Range("E45").Value = "'  +0.36%  "
$ws.Range("D46").Value = "'67.50"
$ws.Range("E46").Value = "'  +1.60%  "
$ws.Range("D47").Value = "'7.339"
$ws.Range("E47").Value = "'  +9.72%  "
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.191"
$ws.Range("E48").Value = "'  +0.01%  "
$ws.Range("D49").Value = "'0.4037"
$ws.Range("E49").Value = "'  -0.75%  "
$ws.Range("B50").Value = "'BabyDogeCoin"
$ws.Range("C50").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000117"
$ws.Range("E50").Value = "'  -1.76%  "
$ws.Range("D51").Value = "'1.703"
$ws.Range("E51").Value = "'  +3.21%  "
